# Apply cryptocurrency price/volume updates as described in the commit:
# "Updated symbol list on Wed Jan 18 17:53:42 UTC 2023 with GitHub Actions"
#
# The source cells are plain text (inline strings) holding formatted
# numeric/percentage values (e.g. "293.86", "-2.66%"). We force each target
# cell to Text format before writing so Excel does not silently reinterpret
# the string as a number (which would normalize things like trailing zeros,
# e.g. "0.0001250" -> "0.000125"). After writing, the style is reset back to
# "Normal" so no stray per-cell formatting is left behind.

function Set-CellText($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" '293.86'
Set-CellText $ws "E2" '-2.66%'
Set-CellText $ws "D3" '31.29'
Set-CellText $ws "E3" '-1.55%'
Set-CellText $ws "D4" '4.972'
Set-CellText $ws "E4" '-0.98%'
Set-CellText $ws "D5" '0.07371'
Set-CellText $ws "E5" '-5.59%'
Set-CellText $ws "D6" '1.843'
Set-CellText $ws "E6" '-9.50%'
Set-CellText $ws "D7" '7.674'
Set-CellText $ws "E7" '-1.42%'
Set-CellText $ws "D8" '3.755'
Set-CellText $ws "E8" '-0.57%'
Set-CellText $ws "D9" '0.9080'
Set-CellText $ws "E9" '-1.00%'
Set-CellText $ws "D10" '0.1648'
Set-CellText $ws "E10" '-5.84%'
Set-CellText $ws "D11" '0.07591'
Set-CellText $ws "E11" '-3.37%'
Set-CellText $ws "D12" '0.08168'
Set-CellText $ws "E12" '-6.97%'
Set-CellText $ws "D13" '0.02989'
Set-CellText $ws "E13" '-4.31%'
Set-CellText $ws "D14" '0.09993'
Set-CellText $ws "E14" '-0.03%'
Set-CellText $ws "D15" '0.001497'
Set-CellText $ws "E15" '-1.09%'
Set-CellText $ws "D16" '0.005656'
Set-CellText $ws "E16" '-2.63%'
Set-CellText $ws "D18" '3.457'
Set-CellText $ws "E18" '-0.17%'
Set-CellText $ws "D19" '2.122'
Set-CellText $ws "E19" '-6.42%'
Set-CellText $ws "E20" '0.38%'
Set-CellText $ws "D21" '0.1315'
Set-CellText $ws "E21" '2.62%'
Set-CellText $ws "D22" '4.327'
Set-CellText $ws "E22" '4.32%'
Set-CellText $ws "D23" '0.1978'
Set-CellText $ws "E23" '9.28%'
Set-CellText $ws "D24" '0.04481'
Set-CellText $ws "E24" '-2.44%'
Set-CellText $ws "D25" '0.001225'
Set-CellText $ws "E25" '-1.13%'
Set-CellText $ws "D26" '0.004046'
Set-CellText $ws "E26" '-9.64%'
Set-CellText $ws "D27" '0.0001250'
Set-CellText $ws "E27" '0.07%'
Set-CellText $ws "D39" '0.01637'
Set-CellText $ws "E39" '-5.52%'
Set-CellText $ws "D40" '0.04404'
Set-CellText $ws "E40" '-7.23%'
Set-CellText $ws "D41" '0.007374'
Set-CellText $ws "E41" '-0.58%'
Set-CellText $ws "D42" '0.1324'
Set-CellText $ws "E42" '-2.34%'
Set-CellText $ws "D43" '0.002055'
Set-CellText $ws "E43" '-3.91%'
Set-CellText $ws "D44" '0.01112'
Set-CellText $ws "E44" '2.86%'
Set-CellText $ws "D45" '0.00005976'
Set-CellText $ws "E45" '-1.49%'
Set-CellText $ws "D46" '0.00000000750'
Set-CellText $ws "E46" '0.07%'
Set-CellText $ws "D47" '1.942'
Set-CellText $ws "E47" '136.77%'
Set-CellText $ws "D48" '0.003005'
Set-CellText $ws "E48" '-14.10%'
Set-CellText $ws "D49" '0.00002100'
Set-CellText $ws "E49" '0.07%'
Set-CellText $ws "D50" '0.0002000'
Set-CellText $ws "E50" '0.07%'
